# daily auto push: 2026-01-21 09:44 UTC
# Insert a new data row for 2026/01/21 (水) just above the existing
# "2026/12/29" block, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 670 (shifts old rows 670..711 down to 671..712)
$ws.Rows("670:670").Insert()

# Populate the newly inserted row with the new record.
# The leading apostrophe forces the date-looking text to stay literal text
# (matching how the rest of column A is stored) instead of Excel
# auto-converting it to a date serial number.
$ws.Range("A670").Value = "'2026/01/21"
$ws.Range("B670").Value = "水"
$ws.Range("C670").Value = 16
$ws.Range("D670").Value = 26

# Undo the "quote-prefixed text" number format that the apostrophe trick
# applies, so the cell matches the plain (unstyled) look of its neighbors.
$ws.Range("A670").Style = "Normal"
